$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The checkpoint list loses its "#6" entry (row 12, column C) -- that day no
# longer has a checkpoint due -- and every later checkpoint (#7..#12) shifts
# down one number to fill the gap (so the final "Checkpoint #12 due" string
# disappears entirely).
$ws.Range("C25").Value = "Checkpoint #11 due"
$ws.Range("C24").Value = "Checkpoint #10 due"
$ws.Range("C22").Value = "Checkpoint #9 due"
$ws.Range("C21").Value = "Checkpoint #8 due"
$ws.Range("C19").Value = "Checkpoint #7 due"
$ws.Range("C13").Value = "Checkpoint #6 due"

# Row 12 (previously "Checkpoint #6 due") becomes a quote-prefixed text note.
$ws.Range("C12").Value = "'**No checkpoint**"

# Reflect the new focus in the saved selection / scroll state.
$ws.Range("C12:C13").Select()
